# Swap the match-data (columns B:AD) between pairs of rows while leaving
# column A (the running index/rank) untouched. This mirrors two rows'
# worth of betting-odds data being re-sorted against each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(131, 132),
    @(163, 164),
    @(183, 184),
    @(192, 193)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
